# Clean up the first sheet of the demo workbook:
# remove the row for "6" (name 6 / location 6), shifting all subsequent
# rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lapa1")
$ws.Rows.Item(6).Delete()
